$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: header values - B1:E1 should carry the same style as A1 (bold/border/center)
$ws.Range("B1:E1").Value = 0
$ws.Range("A1").Copy()
$ws.Range("B1:E1").PasteSpecial(-4122)

$ws.Range("A1").Value = 0
$ws.Range("B1").Value = 0.1
$ws.Range("C1").Value = 0.2
$ws.Range("D1").Value = 0.5
$ws.Range("E1").Value = 1

# Row 2
$ws.Range("A2").Value = 6275.982952620485
$ws.Range("B2").Value = 6301.83023095386
$ws.Range("C2").Value = 6363.372900488996
$ws.Range("D2").Value = 6530.920692923298
$ws.Range("E2").Value = 7243.782941570591

# Row 3
$ws.Range("A3").Value = 11.41064167862783
$ws.Range("B3").Value = 11.4924972835958
$ws.Range("C3").Value = 11.84261183891548
$ws.Range("D3").Value = 11.64127917862189
$ws.Range("E3").Value = 12.266748390855

# Row 4
$ws.Range("A4").Value = 250
$ws.Range("B4").Value = 250
$ws.Range("C4").Value = 250
$ws.Range("D4").Value = 250
$ws.Range("E4").Value = 250

# Row 5
$ws.Range("A5").Value = 0.3936485081193504
$ws.Range("B5").Value = 0.3917046636269098
$ws.Range("C5").Value = 0.3932485401167907
$ws.Range("D5").Value = 0.3872490200783937
$ws.Range("E5").Value = 0.3768498520118391

# Row 6
$ws.Range("A6").Value = 0.0007870296236989592
$ws.Range("B6").Value = 0.0007833554177148746
$ws.Range("C6").Value = 0.0007864257947660684
$ws.Range("D6").Value = 0.0007739510642027224
$ws.Range("E6").Value = 0.0007532670255558446
